# GradientAnalysis.xlsx update
# - Add "Const Values" sheet with a small reference table of TF op constants.
# - Add two new rows to the "Operations" sheet describing the Softmax_grad subgraph,
#   and mark the "Op/tf.name_scope" header cells with the "Good" (green) cell style.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Operations sheet: new rows + "Good" style highlight on the two section
#    header cells (A7 "Neg, Neg_grad" scope & A10 "mul_grad subgraph" scope).
# ---------------------------------------------------------------------------
$opsWs = $wb.Worksheets.Item("Operations")

$opsWs.Range("A7").Style = "Good"
$opsWs.Range("A10").Style = "Good"

$opsWs.Range("A12").Value = "Mul"
$opsWs.Range("B12").Value = "gradients/Log_grad/Reciprocal 1x10" + [char]10 + "gradients/mul_grad/tuple/control_dependency 1x10"
$opsWs.Range("B12").WrapText = $true
$opsWs.Range("C12").Value = "gradients/Softmax_grad/mul 1x10" + [char]10 + "gradients/Softmax_grad/sub 1x10"
$opsWs.Range("C12").WrapText = $true
$opsWs.Rows.Item(12).RowHeight = 29

$opsWs.Range("A13").Value = "Softmax_gr"
$opsWs.Range("A13").Style = "Good"

Write-Output "operations sheet updated"

# ---------------------------------------------------------------------------
# 2. New "Const Values" worksheet, placed after "Operations".
# ---------------------------------------------------------------------------
$cvWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $opsWs)
$cvWs.Name = "Const Values"

$cvWs.Range("A1").Value = "Const Values"
$cvWs.Range("B1").Value = "Value"
$cvWs.Range("C1").Value = "Op"
$cvWs.Range("D1").Value = "Impact"

$cvWs.Range("A4").Value = "gradients/Shape"
$cvWs.Range("B4").Value = "empty"

$cvWs.Range("A5").Value = "gradients/grad_ys_0"
$cvWs.Range("B5").Value = "1 float"

$cvWs.Range("C6").Value = "Fill"
$cvWs.Range("D6").Value = "This operation creates a tensor of shape dims and fills it with value"
$cvWs.Range("D6").WrapText = $true
$cvWs.Rows.Item(6).RowHeight = 58

$cvWs.Range("C7").Value = "Neg"
$cvWs.Range("D7").Value = "Computes numerical negative value element-wise."
$cvWs.Range("D7").WrapText = $true
$cvWs.Rows.Item(7).RowHeight = 58

$cvWs.Range("A8").Value = "gradients/Sum_grad/Reshape/shape"
$cvWs.Range("B8").Value = "[1 1]"
$cvWs.Range("C8").Value = "Const"

$cvWs.Range("A9").Value = "gradients/Sum_grad/Reshape"
$cvWs.Range("C9").Value = "Reshape"
$cvWs.Range("D9").Value = "Here just set the shape 1 1 which is actual a shape of a scalar, no practical change."
$cvWs.Range("D9").WrapText = $true

# E9 is a rich-text cell: "Given tensor, this operation returns a tensor that
# has the same values as tensor with shape shape." with the parameter words
# "tensor"/"tensor"/"shape" typeset in a monospace font.
$e9 = $cvWs.Range("E9")
$e9.Value = "Given tensor, this operation returns a tensor that has the same values as tensor with shape shape."
$e9.WrapText = $true
$e9.Font.Size = 8
$e9.Font.Color = 2367776
$e9.Font.Name = "Noto Sans"
$e9.Characters(7, 6).Font.Name = "Roboto Mono"
$e9.Characters(7, 6).Font.Size = 7
$e9.Characters(7, 6).Font.Color = 5195575
$e9.Characters(75, 6).Font.Name = "Roboto Mono"
$e9.Characters(75, 6).Font.Size = 7
$e9.Characters(75, 6).Font.Color = 5195575
$e9.Characters(93, 5).Font.Name = "Roboto Mono"
$e9.Characters(93, 5).Font.Size = 7
$e9.Characters(93, 5).Font.Color = 5195575
$cvWs.Rows.Item(9).RowHeight = 144

$cvWs.Range("A10").Value = "gradients/Sum_grad/Const"
$cvWs.Range("B10").Value = "[1 10]"
$cvWs.Range("C10").Value = "Const"
$cvWs.Range("D10").Value = "is the input to tile"
$cvWs.Range("D10").WrapText = $true

$cvWs.Range("A11").Value = "gradients/Sum_grad/Tile"
$cvWs.Range("B11").Value = "return a tensor 1x10"

$cvWs.Range("A12").Value = "gradients/mul_grad/Mul"
$cvWs.Range("B12").Value = "multiplies log and tile"

$cvWs.Range("A13").Value = "gradients/mul_grad/Mul1"
$cvWs.Range("B13").Value = "multiplies placeholder and tile"

# ColumnWidth goes through the engine's pixel-quantized character-width model,
# so the inputs below are chosen to land the resulting OOXML width as close as
# possible to the authored widths (31.82 / 28.18 / 12.27 / 17.45 chars).
$cvWs.Columns.Item(1).ColumnWidth = 30.91891
$cvWs.Columns.Item(2).ColumnWidth = 27.25219
$cvWs.Columns.Item(3).ColumnWidth = 11.41703
$cvWs.Columns.Item(4).ColumnWidth = 16.58562

$cvWs.PageSetup.PaperSize = 9
$cvWs.PageSetup.Orientation = 1

Write-Output "const values sheet created"

# ---------------------------------------------------------------------------
# 3. Selections: make "Const Values" the active/selected tab with B13
#    selected, matching where the author was last working.
# ---------------------------------------------------------------------------
$opsWs.Range("B13").Select()
$cvWs.Activate()
$cvWs.Range("B13").Select()

Write-Output "done"
